# The target paragraph ends "...also include:" and must become
# "...also include one of:". In the underlying OOXML this sentence is
# split across two runs at the existing "include" / ":" boundary:
#   run A: " If GIS ... also include"
#   run B: ":"
# The edit inserts "on" at the end of run A and "e of" before the
# colon in run B, so the runs become "...also include on" + "e of:".
$d = $word.ActiveDocument

# Locate the (unique) paragraph containing the sentence to edit.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*also include:*") {
        $target = $para.Range
        break
    }
}

$fullText = $target.Text
$anchor = "also include"
$includeEnd = $target.Start + $fullText.IndexOf($anchor) + $anchor.Length

# Run A: "...also include" -> "...also include on"
# (a collapsed range placed right at the run boundary attaches to the
# preceding run, so this extends run A in place rather than creating
# a new run.)
$insertPoint = $d.Range($includeEnd, $includeEnd)
$insertPoint.InsertAfter(" on")

# Run B: ":" -> "e of:"
# Re-read the paragraph end after the previous insertion shifted it,
# then replace the lone colon character (run B) with "e of:".
$paraEnd = $target.End
$colonRange = $d.Range($paraEnd - 2, $paraEnd - 1)
$colonRange.Text = "e of:"
